# Re-running university responses/police actions analyses
# Updates both summary_counts (sheet1) and response_action_counts (sheet2) sheets
# with refreshed counts and recomputed percentages (count / new total).

$wb = $excel.ActiveWorkbook
$wsSummary = $wb.Worksheets.Item("summary_counts")
$wsActions = $wb.Worksheets.Item("response_action_counts")

# --- summary_counts sheet: updated counts in column B ---
$wsSummary.Range("B2").Value = 868
$wsSummary.Range("B3").Value = 409
$wsSummary.Range("B5").Value = 904
$wsSummary.Range("B6").Value = 686
$wsSummary.Range("B7").Value = 986
$wsSummary.Range("B8").Value = 692
$wsSummary.Range("B9").Value = 816
$wsSummary.Range("B10").Value = 1570
$wsSummary.Range("B11").Value = 911
$wsSummary.Range("B12").Value = 2125
$wsSummary.Range("B13").Value = 3282
$wsSummary.Range("B14").Value = 5407

# --- response_action_counts sheet: updated counts (D) and percentages (E) ---
$wsActions.Range("E2").Value = 0.009617162936933604
$wsActions.Range("E3").Value = 0.000369890882189754
$wsActions.Range("E4").Value = 0.02959127057518032
$wsActions.Range("E5").Value = 0.004438690586277048
$wsActions.Range("D6").Value = 162
$wsActions.Range("E6").Value = 0.02996116145737008
$wsActions.Range("E7").Value = 0.002589236175328278
$wsActions.Range("E8").Value = 0.008877381172554096
$wsActions.Range("E9").Value = 0.0009247272054743851
$wsActions.Range("E10").Value = 0.01054189014240799
$wsActions.Range("E11").Value = 0.004993526909561679
$wsActions.Range("E12").Value = 0.01035694470131311
$wsActions.Range("E13").Value = 0.005363417791751433
$wsActions.Range("E14").Value = 0.01904938043277233
$wsActions.Range("D15").Value = 180
$wsActions.Range("E15").Value = 0.03329017939707786
$wsActions.Range("E16").Value = 0.07601257628999446
$wsActions.Range("E17").Value = 0.001109672646569262
$wsActions.Range("D18").Value = 25
$wsActions.Range("E18").Value = 0.004623636027371925
$wsActions.Range("D19").Value = 65
$wsActions.Range("E19").Value = 0.01202145367116701
$wsActions.Range("E20").Value = 0.005918254115036064
$wsActions.Range("E21").Value = 0.003144072498612909
$wsActions.Range("D22").Value = 89
$wsActions.Range("E22").Value = 0.01646014425744406
$wsActions.Range("E23").Value = 0.05326428703532458
$wsActions.Range("E24").Value = 0.009802108378028482
$wsActions.Range("E25").Value = 0.01461068984649528
$wsActions.Range("E26").Value = 0.01257628999445164
$wsActions.Range("D27").Value = 301
$wsActions.Range("E27").Value = 0.05566857776955798
$wsActions.Range("E28").Value = 0.02441279822452377
$wsActions.Range("D29").Value = 28
$wsActions.Range("E29").Value = 0.005178472350656557
$wsActions.Range("D30").Value = 3
$wsActions.Range("E30").Value = 0.000554836323284631
$wsActions.Range("D31").Value = 461
$wsActions.Range("E31").Value = 0.08525984834473831
$wsActions.Range("E32").Value = 0.02700203439985204
$wsActions.Range("D33").Value = 384
$wsActions.Range("E33").Value = 0.07101904938043277
$wsActions.Range("E34").Value = 0.004253745145182172
$wsActions.Range("D35").Value = 43
$wsActions.Range("E35").Value = 0.007952653967079712
$wsActions.Range("E36").Value = 0.004808581468466802
$wsActions.Range("E37").Value = 0.001294618087664139
$wsActions.Range("D38").Value = 89
$wsActions.Range("E38").Value = 0.01646014425744406
$wsActions.Range("D39").Value = 34
$wsActions.Range("E39").Value = 0.006288144997225818
$wsActions.Range("E40").Value = 0.009987053819123359
$wsActions.Range("D41").Value = 183
$wsActions.Range("E41").Value = 0.03384501572036249
$wsActions.Range("D42").Value = 9
$wsActions.Range("E42").Value = 0.001664508969853893
$wsActions.Range("D43").Value = 4
$wsActions.Range("E43").Value = 0.0007397817643795081
$wsActions.Range("E44").Value = 0.002219345293138524
$wsActions.Range("E45").Value = 0.01683003513963381
$wsActions.Range("D46").Value = 20
$wsActions.Range("E46").Value = 0.00369890882189754
$wsActions.Range("E47").Value = 0.01775476234510819
$wsActions.Range("E48").Value = 0.0136859626410209
$wsActions.Range("D49").Value = 58
$wsActions.Range("E49").Value = 0.01072683558350287
$wsActions.Range("D50").Value = 27
$wsActions.Range("E50").Value = 0.004993526909561679
$wsActions.Range("E51").Value = 0.008507490290364343
$wsActions.Range("D52").Value = 484
$wsActions.Range("E52").Value = 0.08951359348992047
$wsActions.Range("D53").Value = 320
$wsActions.Range("E53").Value = 0.05918254115036065
$wsActions.Range("D54").Value = 331
$wsActions.Range("E54").Value = 0.06121694100240429
$wsActions.Range("D55").Value = 27
$wsActions.Range("E55").Value = 0.004993526909561679
$wsActions.Range("D56").Value = 30
$wsActions.Range("E56").Value = 0.005548363232846311
$wsActions.Range("D57").Value = 186
$wsActions.Range("E57").Value = 0.03439985204364712
$wsActions.Range("E58").Value = 0.000554836323284631
$wsActions.Range("D59").Value = 133
$wsActions.Range("E59").Value = 0.02459774366561864
$wsActions.Range("D60").Value = 35
$wsActions.Range("E60").Value = 0.006473090438320696
$wsActions.Range("D61").Value = 83
$wsActions.Range("E61").Value = 0.01535047161087479
$wsActions.Range("D62").Value = 181
$wsActions.Range("E62").Value = 0.03347512483817274
$wsActions.Range("E63").Value = 0.01683003513963381
$wsActions.Range("E64").Value = 0.002404290734233401
$wsActions.Range("D65").Value = 80
$wsActions.Range("E65").Value = 0.01479563528759016
$wsActions.Range("D66").Value = 196
$wsActions.Range("E66").Value = 0.0362493064545959
$wsActions.Range("D67").Value = 69
$wsActions.Range("E67").Value = 0.01276123543554651
$wsActions.Range("E68").Value = 0.004438690586277048
$wsActions.Range("D69").Value = 180
$wsActions.Range("E69").Value = 0.03329017939707786
$wsActions.Range("E70").Value = 0.001664508969853893
$wsActions.Range("D71").Value = 183
$wsActions.Range("E71").Value = 0.03384501572036249
$wsActions.Range("E72").Value = 0.02145367116700573
$wsActions.Range("D73").Value = 186
$wsActions.Range("E73").Value = 0.03439985204364712
$wsActions.Range("D74").Value = 285
$wsActions.Range("E74").Value = 0.05270945071203995
$wsActions.Range("E75").Value = 0.005548363232846311
$wsActions.Range("D76").Value = 115
$wsActions.Range("E76").Value = 0.02126872572591086
$wsActions.Range("E77").Value = 0.009062326613648974
$wsActions.Range("E78").Value = 0.000554836323284631
$wsActions.Range("D79").Value = 22
$wsActions.Range("E79").Value = 0.004068799704087294
$wsActions.Range("D80").Value = 32
$wsActions.Range("E80").Value = 0.005918254115036064
$wsActions.Range("D81").Value = 5407

Write-Host "Applied updated police/university response counts and recomputed percentages."
